$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from H1 header cell onto new I1/J1 headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

$values = @{
    2  = @(9, 9)
    3  = @(6, 6)
    4  = @(8, 8)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(8, 8)
    8  = @(7, 7)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(9, 9)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
